# test case 5 and 6
# Updates statistics cells on Sheet3, Sheet4, Sheet5 and Sheet7 to reflect
# the newly computed numbers for test case 5 and 6.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet3 ("Statistics")
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

# Row 3 - Number of matched families (AI:AN)
$ws3.Range("AI3").Value = 61
$ws3.Range("AJ3").Value = 47
$ws3.Range("AK3").Value = 43
$ws3.Range("AL3").Value = 58
$ws3.Range("AM3").Value = 60
$ws3.Range("AN3").Value = 74

# Row 4 - Number of matched people (AI:AN)
$ws3.Range("AI4").Value = 143
$ws3.Range("AJ4").Value = 128
$ws3.Range("AK4").Value = 172
$ws3.Range("AL4").Value = 189
$ws3.Range("AM4").Value = 145
$ws3.Range("AN4").Value = 179

# Row 6 - Size 1 (AI:AN), AK6 and AN6 unchanged
$ws3.Range("AI6").Value = 9
$ws3.Range("AJ6").Value = 8
$ws3.Range("AL6").Value = 4
$ws3.Range("AM6").Value = 7

# Row 7 - Size 2 (AI:AN), AK7 unchanged
$ws3.Range("AI7").Value = 29
$ws3.Range("AJ7").Value = 7
$ws3.Range("AL7").Value = 10
$ws3.Range("AM7").Value = 28
$ws3.Range("AN7").Value = 39

# Row 8 - Size 3 (AI:AN)
$ws3.Range("AI8").Value = 17
$ws3.Range("AJ8").Value = 25
$ws3.Range("AK8").Value = 9
$ws3.Range("AL8").Value = 17
$ws3.Range("AM8").Value = 19
$ws3.Range("AN8").Value = 33

# Row 9 - Size 4 (AI:AN), AN9 unchanged
$ws3.Range("AI9").Value = 5
$ws3.Range("AJ9").Value = 4
$ws3.Range("AK9").Value = 25
$ws3.Range("AL9").Value = 21
$ws3.Range("AM9").Value = 5

# Row 10 - Size 5 (AI:AN), AN10 unchanged
$ws3.Range("AI10").Value = 1
$ws3.Range("AJ10").Value = 3
$ws3.Range("AK10").Value = 9
$ws3.Range("AL10").Value = 6
$ws3.Range("AM10").Value = 1

# Scarf column (J) throughout the "Average" section
$ws3.Range("J14").Value = 57
$ws3.Range("J15").Value = 159
$ws3.Range("J16").Value = 5
$ws3.Range("J17").Value = 19
$ws3.Range("J18").Value = 20
$ws3.Range("J19").Value = 10
$ws3.Range("J20").Value = 3

$ws3.Range("J23").Value = 57
$ws3.Range("J24").Value = 143

$ws3.Range("J32").Value = 84
$ws3.Range("J33").Value = 436

# ---------------------------------------------------------------------
# Sheet4 ("Best matched")
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")

$ws4.Range("J19").Value = 67
$ws4.Range("J20").Value = 64
$ws4.Range("J21").Value = 47
$ws4.Range("J22").Value = 66
$ws4.Range("J23").Value = 41
$ws4.Range("J24").Value = 58

$ws4.Range("J28").Value = 184
$ws4.Range("J29").Value = 191
$ws4.Range("J30").Value = 128
$ws4.Range("J31").Value = 150
$ws4.Range("J32").Value = 100
$ws4.Range("J33").Value = 203

# ---------------------------------------------------------------------
# Sheet5 ("Num of families get i-th preferred bundle")
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")

$ws5.Range("N7").Value = 22
$ws5.Range("O7").Value = 15
$ws5.Range("R7").Value = 7
$ws5.Range("S7").Value = 6
$ws5.Range("T7").Value = 26
$ws5.Range("V7").Value = 6
$ws5.Range("W7").Value = 13
$ws5.Range("X7").Value = 6
$ws5.Range("Z7").Value = 1
$ws5.Range("AB7").Value = 15
$ws5.Range("AC7").Value = 1
$ws5.Range("AE7").Value = 7
$ws5.Range("AF7").Value = 12
$ws5.Range("AG7").Value = 6
$ws5.Range("AH7").Value = 17
$ws5.Range("AI7").Value = 14
$ws5.Range("AJ7").Value = 0
$ws5.Range("AK7").Value = 15
$ws5.Range("AM7").Value = 11
$ws5.Range("AP7").Value = 0
$ws5.Range("AQ7").Value = -135

# ---------------------------------------------------------------------
# Sheet7 ("Average bundle rank")
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Sheet7")
$ws7.Range("B7").Value = 24
